$wb = $excel.ActiveWorkbook

# Target OOXML column widths from the handback report regen are
# 29.9777050018311 and 13.7470531463623 characters. Excel's ColumnWidth
# setter only offers 1/6-character (pixel) granularity, so the nearest
# achievable stored widths are 30 and 13.666666666666666 respectively;
# these ColumnWidth inputs (bucket midpoints) land on those values.
$wideColWidth = 29.166666666666664
$narrowColWidth = 12.833333333333332

# --- Overview sheet ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "Handed back: in sync with en-US"
$ws1.Range("F2").Value = "Handed back: in sync with en-US"
$ws1.Range("E3").Value = "Handed back: in sync with en-US"
$ws1.Range("F3").Value = "Handed back: in sync with en-US"
$ws1.Columns.Item(5).ColumnWidth = $wideColWidth
$ws1.Columns.Item(6).ColumnWidth = $wideColWidth

# --- zh-cn sheet ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("K2").Value = "2016-11-08 22:47:33"
$ws2.Range("K3").Value = "2016-11-08 22:47:33"
$ws2.Range("P2").Value = ""
$ws2.Range("P3").Value = ""
$ws2.Columns.Item(3).ColumnWidth = $wideColWidth
$ws2.Columns.Item(16).ColumnWidth = $narrowColWidth

# --- de-de sheet ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("K2").Value = "2016-11-08 22:47:51"
$ws3.Range("K3").Value = "2016-11-08 22:47:51"
$ws3.Range("P2").Value = ""
$ws3.Range("P3").Value = ""
$ws3.Columns.Item(3).ColumnWidth = $wideColWidth
$ws3.Columns.Item(16).ColumnWidth = $narrowColWidth
